$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing data
# (symboling .. gas) one column to the right (B:T -> C:U).
$ws.Columns("B:B").Insert()

# The insert copies column A's (bold/bordered) formatting into the new
# column; the data cells (B2:B9) should stay unformatted like the rest
# of the numeric cells, so strip that back off.
$ws.Range("B2:B9").ClearFormats()

# New column B holds the "Unnamed: 0" index-column summary stats that
# pandas' describe() produced once the DataFrame's row index was
# included as a plain column.
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("B2").Value = 201
$ws.Range("B3").Value = 100
$ws.Range("B4").Value = 58.16786054171152
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 50
$ws.Range("B7").Value = 100
$ws.Range("B8").Value = 150
$ws.Range("B9").Value = 200

# Match the bold/bordered header formatting used by the other header
# cells in row 1.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
